{"js": "// The commit turns the \"18.02.\" date line and the \"57. Methods and\n// Functions Homework Overview\" heading line into bold, 14pt (sz/szCs\n// 28 half-points) heading text, matching the formatting convention\n// used by the other day/lesson headings earlier in this document\n// (<w:b/> instead of <w:bCs/>, sz/szCs bumped from 24 -> 28).\n//\n// Plain Font.bold / Font.size assignment only ever *adds* <w:b/> and\n// leaves the existing <w:bCs/> in place (or turns it into an explicit\n// <w:bCs w:val=\"0\"/>) instead of removing it, so we splice in the\n// exact target run/paragraph-mark formatting via insertOoxml(...,\n// Word.InsertLocation.replace) instead - this lets us reproduce the\n// target markup (<w:b/>, <w:sz w:val=\"28\"/>, <w:szCs w:val=\"28\"/>,\n// with no stray <w:bCs/>) precisely, while keeping every paragraph's\n// identity (w14:paraId/w14:textId/rsid*) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst W_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"';\n\nfunction wrapPackage(paragraphXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document ' + W_NS + '><w:body>' + paragraphXml + '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\n// Locate the two target paragraphs by their text rather than a\n// hard-coded index, so the script is resilient to minor shifts.\nlet dateParaIndex = -1;\nlet headingParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text === \"18.02.\") {\n    dateParaIndex = i;\n  } else if (text === \"57. Methods and Functions Homework Overview\") {\n    headingParaIndex = i;\n  }\n}\n\nif (dateParaIndex === -1 || headingParaIndex === -1) {\n  throw new Error(\"Could not locate the target paragraphs.\");\n}\n\n// \"18.02.\" paragraph: both the paragraph mark run properties (pPr/rPr)\n// and the run's own rPr change from bCs/sz24/szCs24 to b/sz28/szCs28.\nconst dateParagraphXml =\n  '<w:p w14:paraId=\"3D091D97\" w14:textId=\"712FF981\" w:rsidR=\"002B2F0D\" w:rsidRDefault=\"002B2F0D\" w:rsidP=\"00713937\">' +\n  '<w:pPr><w:spacing w:after=\"0\"/><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>18.02.</w:t></w:r>' +\n  '</w:p>';\n\ncontext.document.body.paragraphs.items[dateParaIndex].insertOoxml(\n  wrapPackage(dateParagraphXml),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// \"57. Methods and Functions Homework Overview\" paragraph: only the\n// run's rPr changes; the paragraph mark's rPr (pPr/rPr) is left as-is\n// (bCs/sz24/szCs24).\nconst headingParagraphXml =\n  '<w:p w14:paraId=\"5883B070\" w14:textId=\"756C418F\" w:rsidR=\"002B2F0D\" w:rsidRDefault=\"002B2F0D\" w:rsidP=\"00713937\">' +\n  '<w:pPr><w:spacing w:after=\"0\"/><w:rPr><w:bCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>' +\n  '<w:r w:rsidRPr=\"002B2F0D\"><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>57. Methods and Functions Homework Overview</w:t></w:r>' +\n  '</w:p>';\n\ncontext.document.body.paragraphs.items[headingParaIndex].insertOoxml(\n  wrapPackage(headingParagraphXml),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# The commit turns the \"18.02.\" date line and the \"57. Methods and\n# Functions Homework Overview\" heading line into bold, 14pt (sz/szCs\n# 28 half-points) heading text, matching the formatting convention\n# used by the other day/lesson headings earlier in this document\n# (<w:b/> instead of <w:bCs/>, sz/szCs bumped from 24 -> 28).\n#\n# Plain Font.Bold / Font.Size assignment only ever *adds* <w:b/> and\n# leaves the pre-existing <w:bCs/> in place (or turns it into an\n# explicit <w:bCs w:val=\"0\"/>) instead of removing it, so we splice in\n# the exact target run / paragraph-mark formatting via\n# Range.InsertXML(...) instead - this reproduces the target markup\n# (<w:b/>, <w:sz w:val=\"28\"/>, <w:szCs w:val=\"28\"/>, with no stray\n# <w:bCs/>) precisely, while keeping each paragraph's identity\n# (w14:paraId/w14:textId/rsid*) untouched.\n\n$d = $word.ActiveDocument\n\nfunction New-PackageXml([string]$paragraphXml) {\n    $wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"'\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document ' + $wNs + '><w:body>' + $paragraphXml + '</w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# Locate the two target paragraphs by their text rather than a\n# hard-coded index, so the script is resilient to minor shifts.\n$dateParagraph = $null\n$headingParagraph = $null\nforeach ($para in $d.Paragraphs) {\n    $text = $para.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($text -eq \"18.02.\") {\n        $dateParagraph = $para\n    } elseif ($text -eq \"57. Methods and Functions Homework Overview\") {\n        $headingParagraph = $para\n    }\n}\n\n# \"18.02.\" paragraph: both the paragraph mark run properties (pPr/rPr)\n# and the run's own rPr change from bCs/sz24/szCs24 to b/sz28/szCs28.\n$dateParagraphXml = '<w:p w14:paraId=\"3D091D97\" w14:textId=\"712FF981\" w:rsidR=\"002B2F0D\" w:rsidRDefault=\"002B2F0D\" w:rsidP=\"00713937\">' +\n    '<w:pPr><w:spacing w:after=\"0\"/><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>18.02.</w:t></w:r>' +\n    '</w:p>'\n$dateParagraph.Range.InsertXML((New-PackageXml $dateParagraphXml))\n\n# \"57. Methods and Functions Homework Overview\" paragraph: only the\n# run's rPr changes; the paragraph mark's rPr (pPr/rPr) is left as-is\n# (bCs/sz24/szCs24).\n$headingParagraphXml = '<w:p w14:paraId=\"5883B070\" w14:textId=\"756C418F\" w:rsidR=\"002B2F0D\" w:rsidRDefault=\"002B2F0D\" w:rsidP=\"00713937\">' +\n    '<w:pPr><w:spacing w:after=\"0\"/><w:rPr><w:bCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>' +\n    '<w:r w:rsidRPr=\"002B2F0D\"><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>57. Methods and Functions Homework Overview</w:t></w:r>' +\n    '</w:p>'\n$headingParagraph.Range.InsertXML((New-PackageXml $headingParagraphXml))\n"}
